$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 24 with the new test-mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A24").Value = "Kun jij dit afhandelen?"
$ws.Range("B24").Value = "mailmind.test@zohomail.eu"
$ws.Range("C24").Value = "Testmail #3: Kun jij dit afhandelen?"
$ws.Range("D24").Value = "Planning / Afspraak"
$ws.Range("E24").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Range("F24").Value = "2025-08-05 18:10:55"
$ws.Range("G24").Value = "Ja"
$ws.Range("H24").Value = "Ja"
$ws.Range("I24").Value = "Nee"
$ws.Range("J24").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row 24 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "23")
    $newRange = $ws.Range($col + "2:" + $col + "24")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Sheet "Dashboard": bump the "Planning / Afspraak" count from 17 to 18 ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 18
